# Estadisticos Segundo Parcial 26 Mayo
# Insert a new rescatable student row at the top of the data in "Rescatables",
# pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new row above current row 2 (the first data row), shifting
# the existing data rows down by one. Excel copies the formatting of the
# row above (the bold header row) into the freshly inserted row, so clear
# that back to the default (unstyled) look used by the other data rows.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:G2").ClearFormats()

# Populate the new row with the new rescatable student's data.
$ws.Range("A2").Value = 23330051920313
$ws.Range("B2").Value = "VIVANCO"
$ws.Range("C2").Value = "VIVANCO"
$ws.Range("D2").Value = "LUIS AARON"
$ws.Range("E2").Value = "Ingles IV"
$ws.Range("F2").Value = "4APM"
$ws.Range("G2").Value = 4
